# Remove R1 from the BOM (and the implicit ROM/DK-Order mirror) to address
# #42 and prolly #35.
#
# "R1, R18, R19" (Value=1k, DK=YAG2306CT-ND, DNP) was a single BOM line for
# three reference designators. R1 itself is being removed from the design;
# R18 and R19 remain on the board (still DNP, still 1k) but no longer have
# a Digikey part number since that was R1's line item.
#
# Net effect on the "BOM" sheet:
#   - the old row for "R1, R18, R19" / 1k / YAG2306CT-ND / DNP is deleted
#   - every row below it shifts up by one
#   - a replacement row "R18, R19" / 1k / DNP (no DK/PARTNO) is (re)created
#     right after R17, so the row count / layout below it is unchanged
#
# The "DK Order" sheet never listed R1 (it was DNP, so not orderable) so it
# does not need any row-level edit; Excel keeps its shared-string references
# in sync automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Delete the "R1, R18, R19" row entirely - everything below shifts up.
$ws.Rows(33).Delete()

# Re-insert a row for the leftover "R18, R19" DNP placeholder right after
# R17 (which is now row 43), restoring the original row count/positions for
# everything from the old "R21, R22" row onward.
$ws.Rows(44).Insert()

$ws.Range("A44").Value = 2
$ws.Range("B44").Value = "R18, R19"
$ws.Range("C44").Value = "1k"
$ws.Range("F44").Value = "T"
